$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column F width (COM ColumnWidth excludes the ~0.8333 padding that ends
# up in the raw OOXML "width" attribute, so subtract 5/6 to land on 49.5).
$ws.Columns("F").ColumnWidth = 48.666666666666664

# Apply wrap text to the F column header (bordered header style) and body
# cells (plain style), matching the new cellXfs entries.
$ws.Range("F1").WrapText = $true
$ws.Range("F2:F13").WrapText = $true

# New PMIDs discovered for rows 3 and 7, centered/top aligned like the
# existing PMID cells. Build the style on G3 first (wrap, then horizontal,
# then vertical so the engine folds all three into a single new cellXfs
# entry), then copy/paste that format onto the other PMID cells so they all
# share the same style index instead of each minting their own.
$ws.Range("G3").Value = 38682164
$ws.Range("G7").Value = 23894501

$ws.Range("G3").WrapText = $true
$ws.Range("G3").HorizontalAlignment = -4108
$ws.Range("G3").VerticalAlignment = -4160

$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights, auto-computed by Excel after wrapping the Accession column at
# its new width.
$ws.Rows(1).RowHeight = 16
$ws.Rows(2).RowHeight = 32
$ws.Rows(3).RowHeight = 96
$ws.Rows(4).RowHeight = 16
$ws.Rows(5).RowHeight = 16
$ws.Rows(6).RowHeight = 16
$ws.Rows(7).RowHeight = 16
$ws.Rows(8).RowHeight = 16
$ws.Rows(9).RowHeight = 32
$ws.Rows(10).RowHeight = 112
$ws.Rows(11).RowHeight = 80
$ws.Rows(12).RowHeight = 32
$ws.Rows(13).RowHeight = 96

# Move the selection to F8 (single cell) as in the final saved view.
$ws.Range("F8").Select()
